# Incomes.xlsx - add a new income record (row 93) to the "Incomes" table.
# Commit message: "I've computed the expense average for each category and
# month; computed the frequencies of each category." -> a new income entry
# (Parents transfer, 5 zl, 2025-02-23) was logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Incomes")

$lastRow = 92
$newRow = 93

# Copy the formatting (number formats/styles) of the last data row down to
# the new row first, so the new cells inherit the same date / currency
# styles without introducing new style/numFmt entries.
$ws.Range("A$lastRow`:C$lastRow").Copy() | Out-Null
$ws.Range("A$newRow`:C$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new income record.
$ws.Range("A$newRow").Value2 = 45711          # Income date  -> 2025-02-23
$ws.Range("B$newRow").Value2 = 5              # Income amount
$ws.Range("C$newRow").Value = "Parents transfer"  # Income category

# Grow the Excel table ("Table7") so it covers the newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C$newRow")) | Out-Null

# Reflect the updated viewport/selection like the author had it.
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("C$newRow").Select() | Out-Null
